$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5862395763397217
$ws.Range("B1").Value = 1.250509023666382
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.758764386177063
$ws.Range("E1").Value = 1.499902725219727
